# [fix][cloudkitty] remove filtering and projectId columns
#
# The "Project ID" column (C) is no longer emitted by the report, so drop
# it from the template and shift "Resources" (previously D) left into C.
# Also drop the "(Timezone: UTC)" suffix from the two date-range headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Project ID" column entirely -- this shifts every column to
# its right (D "Resources" -> C, ..., M "\n" -> L) and shrinks the sheet's
# used range from A1:M68 down to A1:L68, matching the diff.
$ws.Columns("C").Delete()

# Rename the two remaining date-range headers.
$ws.Range("A1").Value = "Begin"
$ws.Range("B1").Value = "End"

# Match the new selection recorded in the template: the whole of column C
# (now "Resources") is selected, anchored at C1.
$ws.Columns("C").Select() | Out-Null
